# Add a new "hideInContents" column (O) to the survey sheet and flag the
# two existing "note" rows (row 2 and row 10) as hidden-in-contents (TRUE),
# per the commit: "hide notes in contents screen".

$wb = $excel.ActiveWorkbook
$survey = $wb.Worksheets.Item("survey")

# New header cell + column width to match the sheet's existing header style.
$survey.Range("O1").Value = "hideInContents"
$survey.Columns.Item(15).ColumnWidth = 13.5

# Mark the two "note" rows as hidden from the contents screen.
$survey.Range("O2").Value = $true
$survey.Range("O10").Value = $true

# "settings" becomes the active/selected sheet in the saved workbook.
$wb.Worksheets.Item("settings").Activate()
